$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it currently sits right after
#    the "Features" heading paragraph).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Collapse the two trailing empty paragraphs at the very end of the body
#    into a single empty paragraph (delete one of them).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$pSecondLast = $d.Paragraphs.Item($count - 1)
$pLast = $d.Paragraphs.Item($count)
# Paragraph.Range.Text for an "empty" paragraph still contains the
# paragraph-mark character, so trim before comparing against "".
if ($pSecondLast.Range.Text.Trim() -eq "" -and $pLast.Range.Text.Trim() -eq "") {
    $trailingRange = $d.Range($pSecondLast.Range.Start, $pLast.Range.End)
    $trailingRange.Delete()
}

# ---------------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark right after the "Language Used: C/C++"
#    text, i.e. at the end of that paragraph (before its paragraph mark).
#
#    Directly building a zero-length Range exactly at a paragraph's last
#    character / paragraph-mark boundary confuses Bookmarks.Add in this
#    runtime, so instead: insert a throw-away character there, wrap a
#    (non-collapsed) bookmark range around it, then delete the character.
#    The bookmark collapses down to the correct zero-width location and
#    survives the deletion.
# ---------------------------------------------------------------------------
$targetParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Language Used: C/C++*") {
        $targetParaIndex = $i
    }
}

if ($targetParaIndex -ne -1) {
    $targetPara = $d.Paragraphs.Item($targetParaIndex)
    $insertPos = $targetPara.Range.End - 1

    $insertionPoint = $d.Range($insertPos, $insertPos)
    $insertionPoint.InsertAfter("Z")

    $markerRange = $d.Range($insertPos, $insertPos + 1)
    $d.Bookmarks.Add("_GoBack", $markerRange)

    $markerRange2 = $d.Range($insertPos, $insertPos + 1)
    $markerRange2.Text = ""
}
